# Commit: add OR stunting given Zinc; preliminary intervention coverages
#
# 1. Insert a new worksheet "OR stunting Zinc" right after "OR stunting
#    diarrhoea" (and before "birth distribution"), holding the odds-ratio
#    of stunting given zinc supplementation by age band.
# 2. Append a new worksheet "Intervention coverages" at the very end of
#    the workbook, holding pre-2016 coverage rates for two interventions.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "OR stunting Zinc" — inserted after "OR stunting diarrhoea"
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("OR stunting diarrhoea")
$zincSheet = $wb.Worksheets.Add($null, $afterSheet)
$zincSheet.Name = "OR stunting Zinc"

$zincSheet.Range("A1").Value = "<1 month"
$zincSheet.Range("B1").Value = "1-5 months"
$zincSheet.Range("C1").Value = "6-11 months"
$zincSheet.Range("D1").Value = "12-23 months"
$zincSheet.Range("E1").Value = "24-59 months"

$zincSheet.Range("A2").Value = 0.9
$zincSheet.Range("B2").Value = 0.9
$zincSheet.Range("C2").Value = 0.9
$zincSheet.Range("D2").Value = 0.9
$zincSheet.Range("E2").Value = 0.9

# ---------------------------------------------------------------------------
# 2. "Intervention coverages" — appended as the last sheet
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$covSheet = $wb.Worksheets.Add($null, $lastSheet)
$covSheet.Name = "Intervention coverages"

$covSheet.Columns.Item(1).ColumnWidth = 33.03
$covSheet.Columns.Item(2).ColumnWidth = 14.74

$covSheet.Range("A1").Value = "Intervention"
$covSheet.Range("B1").Value = "pre-2016"

$covSheet.Range("A2").Value = "Zinc supplementation"
$covSheet.Range("B2").Value = 0.4

$covSheet.Range("A3").Value = "Vitamin A supplementation"
$covSheet.Range("B3").Value = 0.5
